# "add guild data module" - append a new Property row describing the
# guild id field (GuilID / object / Friend / 工会ID), mirroring the
# existing rows on the "Property" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$row = 11

# Match the shared-string insertion order baked into the target file:
# J (工会ID), then B (object), then A (GuilID); I reuses the existing
# "Friend" string.
$ws.Cells.Item($row, 10).Value = "工会ID"
$ws.Cells.Item($row, 2).Value = "object"
$ws.Cells.Item($row, 1).Value = "GuilID"

$ws.Cells.Item($row, 3).Value = $true
$ws.Cells.Item($row, 4).Value = $true
$ws.Cells.Item($row, 5).Value = $true
$ws.Cells.Item($row, 6).Value = $true
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = "Friend"

# Text columns (A, B, I, J) carry the "@" text number format used
# throughout the sheet (style index 1).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 10).NumberFormat = "@"

$ws.Range("E19").Select()
